# Fruta / hortaliza, semanal
# Insert a new weekly block of 5 data rows (report date serial 45204) right before the
# existing row 321 ("Cultivar IV Region" / "Cultivar V Region" entries), pushing the
# rest of the Chirimoya price table down by 5 rows (321-412 -> 326-417).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new blank rows starting at row 321 (existing rows 321.. shift down to 326..)
$ws.Rows.Item(321).Resize(5).Insert()

# Common (unchanged across this block) column values
$A = 3
$B = "Femacal de La Calera"
$C = "Coquimbo"
$E = 5
$F = "Fruta"
$G = 100107
$H = "Otros"
$I = 100107002
$J = "Chirimoya"
$Q = "`$/bandeja 10 kilos"
$T = 10

# New report date for this block: serial 45204 (2023-10-05)
$fecha = [DateTime]::FromOADate(45204)

function Set-ChirimoyaRow($Row, $Variedad, $Calidad, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $Origen, $PrecioKg) {
    $ws.Cells.Item($Row, 1).Value = $A
    $ws.Cells.Item($Row, 2).Value = $B
    $ws.Cells.Item($Row, 3).Value = $C
    $ws.Cells.Item($Row, 4).Value = $fecha
    $ws.Cells.Item($Row, 5).Value = $E
    $ws.Cells.Item($Row, 6).Value = $F
    $ws.Cells.Item($Row, 7).Value = $G
    $ws.Cells.Item($Row, 8).Value = $H
    $ws.Cells.Item($Row, 9).Value = $I
    $ws.Cells.Item($Row, 10).Value = $J
    $ws.Cells.Item($Row, 11).Value = $Variedad
    $ws.Cells.Item($Row, 12).Value = $Calidad
    $ws.Cells.Item($Row, 13).Value = $Volumen
    $ws.Cells.Item($Row, 14).Value = $PrecioMin
    $ws.Cells.Item($Row, 15).Value = $PrecioMax
    $ws.Cells.Item($Row, 16).Value = $PrecioProm
    $ws.Cells.Item($Row, 17).Value = $Q
    $ws.Cells.Item($Row, 18).Value = $Origen
    $ws.Cells.Item($Row, 19).Value = $PrecioKg
    $ws.Cells.Item($Row, 20).Value = $T
}

Set-ChirimoyaRow 321 "Cultivar IV Región" "Especial" 67 30000 30000 30000 "Provincia del Elquí" 3000
Set-ChirimoyaRow 322 "Cultivar IV Región" "Primera"  68 28000 28000 28000 "Provincia del Elquí" 2800
Set-ChirimoyaRow 323 "Cultivar IV Región" "Segunda"  56 25000 25000 25000 "Provincia del Elquí" 2500
Set-ChirimoyaRow 324 "Cultivar V Región"  "Primera"  45 27000 27000 27000 "Provincia de Quillota" 2700
Set-ChirimoyaRow 325 "Cultivar V Región"  "Segunda"  45 23000 23000 23000 "Provincia de Quillota" 2300
